$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# --- Row 15: remove its only cell (B15) entirely, which removes the row itself ---
$ws.Range("B15").Clear()

# --- Row 16: B16 and C16 become blank placeholder cells (keep default/no style) ---
$styleA16 = $ws.Range("A16").Style
$ws.Range("B16").ClearContents()
$ws.Range("B16").Style = $styleA16
$ws.Range("C16").ClearContents()
$ws.Range("C16").Style = $styleA16

# --- Row 17: B17 becomes a blank placeholder cell (no style); C17 keeps its
#     quote-prefixed (s=2) style but loses its value ---
$ws.Range("B17").ClearContents()
$ws.Range("B17").Style = $ws.Range("A17").Style
$ws.Range("C17").ClearContents()
$ws.Range("C17").QuotePrefix = $true

# --- Rows 18-21: B cell disappears completely; C cell keeps its s=2 style, blank ---
$ws.Range("B18").Clear()
$ws.Range("C18").ClearContents()
$ws.Range("C18").QuotePrefix = $true

$ws.Range("B19").Clear()
$ws.Range("C19").ClearContents()
$ws.Range("C19").QuotePrefix = $true

$ws.Range("B20").Clear()
$ws.Range("C20").ClearContents()
$ws.Range("C20").QuotePrefix = $true

$ws.Range("B21").Clear()
$ws.Range("C21").ClearContents()
$ws.Range("C21").QuotePrefix = $true

# --- Row 22: B22 disappears completely; C22 keeps its s=2 style, blank ---
$ws.Range("B22").Clear()
$ws.Range("C22").ClearContents()
$ws.Range("C22").QuotePrefix = $true

# --- Update the selection shown on the sheet ---
$ws.Range("A15:I29").Select()

# --- Update the workbook window geometry ---
$win = $excel.Windows.Item(1)
$win.Left = 38290
$win.Top = -110
$win.Width = 38620
$win.Height = 21360
